$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = "MCT-1A-Tecnologia dos Materiais"

# Row 4
$ws.Range("D4").Value = "MEC-1A-Tecnologia dos Materiais"
$ws.Range("F4").Value = "MEC-1A-Tecnologia dos Materiais"

# Row 6
$ws.Range("D6").Value = "MEC-1A-Tecnologia dos Materiais"
$ws.Range("E6").Value = "MCT-1A-Tecnologia dos Materiais"
$ws.Range("F6").Value = "-"

# Row 7
$ws.Range("E7").Value = "MCT-1A-Tecnologia dos Materiais"

# Row 8
$ws.Range("E8").Value = "-"
